$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.051776313236694
$ws.Range("D2").Value = 0.3043183619633165

$ws.Range("C3").Value = -1.209121655653234
$ws.Range("D3").Value = 0.2394506238967216

$ws.Range("C4").Value = -1.335272563459644
$ws.Range("D4").Value = 0.1954398540995492

$ws.Range("C5").Value = -1.168300679555073
$ws.Range("D5").Value = 0.2551905367302996

$ws.Range("C6").Value = -0.6953113875240868
$ws.Range("D6").Value = 0.4941397186039669

$ws.Range("C7").Value = -0.8234208044191134
$ws.Range("D7").Value = 0.4191057862169671

$ws.Range("C8").Value = -0.8687441196471003
$ws.Range("D8").Value = 0.3943659695397073

$ws.Range("C9").Value = 0.1544654603758045
$ws.Range("D9").Value = 0.8786510437444703

$ws.Range("C10").Value = -0.1632521045315588
$ws.Range("D10").Value = 0.8718103136404594

$ws.Range("C11").Value = -0.3176263988856375
$ws.Range("D11").Value = 0.7537633943134667

$wb.Save()
